$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column D: header + values for the week of Feb 19-25, 2024
$ws.Range("D1").Value = "25_02_2024"
$ws.Range("D2").Value = 990
$ws.Range("D3").Value = 934
$ws.Range("D4").Value = 1456
$ws.Range("D5").Value = 2906
$ws.Range("D6").Value = 55

# Update selection to reflect where the user left off (D7)
$ws.Range("D7").Select()
